# Fixed a bug in MultiLevelReplaceReel
# The data rows (reel groups) were reordered. Capture the original values
# for the affected rows first, then write them back out to their new
# destination rows, so we don't rely on manually transcribed numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of destination row -> source row (both refer to the ORIGINAL,
# pre-edit layout of the sheet).
$rowMap = @{
    3  = 7
    4  = 9
    5  = 10
    6  = 11
    7  = 12
    9  = 5
    10 = 4
    11 = 6
    12 = 15
    13 = 3
    14 = 13
    15 = 14
    16 = 17
    17 = 18
    18 = 16
    19 = 20
    20 = 19
}

# Snapshot original values (columns A-F) for every row referenced above,
# before any writes happen.
$snapshot = @{}
foreach ($srcRow in $rowMap.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowValues = @()
        for ($col = 1; $col -le 6; $col++) {
            $rowValues += $ws.Cells.Item($srcRow, $col).Value2
        }
        $snapshot[$srcRow] = $rowValues
    }
}

# Now write the snapshotted values into their destination rows.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $values = $snapshot[$srcRow]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($destRow, $col).Value = $values[$col - 1]
    }
}
